# Add functions to find available rooms and choose desired room
#
# The raw "day" (date) and "hour" (time-of-day fraction) columns are
# collapsed into a single "date_time" timestamp column so a booking slot
# can be looked up/compared as one value. "available_places" shifts left
# to take the vacated column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F held the date (whole-number serial), column G held the time of
# day (fractional serial), column H held available_places.
# Combine F+G into a single date_time value in F, and move H's values
# into G (which becomes the new available_places column).
for ($row = 2; $row -le 15; $row++) {
    $dateVal  = $ws.Cells.Item($row, 6).Value2()
    $timeVal  = $ws.Cells.Item($row, 7).Value2()
    $availVal = $ws.Cells.Item($row, 8).Value2()

    $ws.Cells.Item($row, 6).Value = $dateVal + $timeVal
    $ws.Cells.Item($row, 7).Value = $availVal
}

# Headers: F becomes date_time, G becomes available_places, H is dropped.
$ws.Range("F1").Value = "date_time"
$ws.Range("G1").Value = "available_places"
$ws.Range("H1:H15").ClearContents()

# date_time gets a combined date+time display format; available_places
# goes back to the plain/general look the old column had.
$ws.Range("F2:F15").NumberFormat = "d/m/yy h:mm;@"
$ws.Range("G2:G15").ClearFormats()

# The merged date_time column is wider than the old bestFit "day" column.
$ws.Columns.Item(6).ColumnWidth = 15.5

[void]$ws.Range("F1").Select()
